# "na migratie db versie 1.24"
#
# The SharePoint content-type metadata (customXml/item1.xml, the
# ct:contentTypeSchema part, plus its auto-maintained itemProps1.xml
# companion) was re-synced after a document-library migration:
#   - ma:contentTypeVersion bumped 16 -> 17
#   - ma:versionID / ma:fieldsID refreshed to new GUID-ish hashes
#   - a new "MediaServiceSearchProperties" field (ns2, dms:Note) was
#     added to the schema (both the <xsd:all> ref and its definition)
#   - the itemProps1 datastore item got a fresh ds:itemID and lost its
#     explicit schemaRefs list (Word regenerates that automatically)
#
# This is done the way Word's object model actually exposes editing of
# package-level custom XML parts: CustomXMLPart.XML is read-only, so to
# change one you locate it (by namespace), read its .XML, transform the
# text, delete the old part and Add() the replacement. Word re-derives
# the companion itemProps part (new ds:itemID, refreshed schemaRefs)
# automatically when it does this.

$d = $word.ActiveDocument

$contentTypeNs  = "http://schemas.microsoft.com/office/2006/metadata/contentType"
$newVersionID   = "203f714213f9939ffdb9d12666e4ed24"
$newFieldsID    = "1026762a0c59493e2a3d7946b4380e53"

$parts = $d.CustomXMLParts
$target = $null

# Find the contentTypeSchema part (there's exactly one such part in this
# document) by probing the parts for the ct: namespace / known marker text.
for ($i = 1; $i -le $parts.Count; $i++) {
    $p = $parts.Item($i)
    if ($p.XML -like "*$contentTypeNs*" -and $p.XML -like "*contentTypeSchema*") {
        $target = $p
    }
}
if ($target -eq $null) {
    $byNs = $parts.SelectByNamespace($contentTypeNs)
    if ($byNs -ne $null -and $byNs.Count -ge 1) {
        $target = $byNs.Item(1)
    }
}

if ($target -ne $null) {
    $xml = $target.XML

    $xml = $xml.Replace(
        'ma:contentTypeVersion="16" ma:contentTypeDescription="Een nieuw document maken." ma:contentTypeScope="" ma:versionID="6f76822a87dec94ded47590863c36fbe"',
        'ma:contentTypeVersion="17" ma:contentTypeDescription="Een nieuw document maken." ma:contentTypeScope="" ma:versionID="' + $newVersionID + '"'
    )

    $xml = $xml.Replace(
        'ma:fieldsID="2b111738bb78e2587f2d436f84d12744"',
        'ma:fieldsID="' + $newFieldsID + '"'
    )

    $xml = $xml.Replace(
        '<xsd:element ref="ns2:MediaServiceObjectDetectorVersions" minOccurs="0"/></xsd:all>',
        '<xsd:element ref="ns2:MediaServiceObjectDetectorVersions" minOccurs="0"/><xsd:element ref="ns2:MediaServiceSearchProperties" minOccurs="0"/></xsd:all>'
    )

    $oldDef = '<xsd:element name="MediaServiceObjectDetectorVersions" ma:index="23" nillable="true" ma:displayName="MediaServiceObjectDetectorVersions" ma:hidden="true" ma:indexed="true" ma:internalName="MediaServiceObjectDetectorVersions" ma:readOnly="true"><xsd:simpleType><xsd:restriction base="dms:Text"/></xsd:simpleType></xsd:element></xsd:schema>'
    $newDef = '<xsd:element name="MediaServiceObjectDetectorVersions" ma:index="23" nillable="true" ma:displayName="MediaServiceObjectDetectorVersions" ma:hidden="true" ma:indexed="true" ma:internalName="MediaServiceObjectDetectorVersions" ma:readOnly="true"><xsd:simpleType><xsd:restriction base="dms:Text"/></xsd:simpleType></xsd:element><xsd:element name="MediaServiceSearchProperties" ma:index="24" nillable="true" ma:displayName="MediaServiceSearchProperties" ma:hidden="true" ma:internalName="MediaServiceSearchProperties" ma:readOnly="true"><xsd:simpleType><xsd:restriction base="dms:Note"/></xsd:simpleType></xsd:element></xsd:schema>'
    $xml = $xml.Replace($oldDef, $newDef)

    $target.Delete()
    $parts.Add($xml) | Out-Null
}
